$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 25 ("Sohyang") keeps its answer but gets a new clip; row 26
# ("Dean") also gets a new clip. These two writes must happen first (in
# this order) so freshly-appended shared strings land at the same table
# slots as the target workbook.
$ws.Range("B26").Value = "LYZ7gdFkmcs"
$ws.Range("B25").Value = "Asmt-dgB-64"
$ws.Range("C25").Value = 68

# --- 16 new quiz rows (ids 38-53, sheet rows 39-54) ---
$ws.Range("D39").Value = "노브레인"
$ws.Range("D40").Value = "노라조"
$ws.Range("B39").Value = "kFv1IQLekU0"
$ws.Range("B40").Value = "ao58vQDMVlQ"
$ws.Range("D41").Value = "허각"
$ws.Range("D42").Value = "김필"
$ws.Range("B41").Value = "3xdfBwFb2DU"
$ws.Range("D43").Value = "어반자카파"
$ws.Range("B42").Value = "T5_O38Bpeto"
$ws.Range("D44").Value = "자우림"
$ws.Range("D45").Value = "체리필터"
$ws.Range("B43").Value = "fmq2k0MkZ0g"
$ws.Range("D46").Value = "김종국"
$ws.Range("D47").Value = "이승기"
$ws.Range("D48").Value = "안예은"
$ws.Range("B44").Value = "qvJ1FHRR1n8"
$ws.Range("B45").Value = "HkQ6HvW9nrA"
$ws.Range("D49").Value = "김상민"
$ws.Range("D50").Value = "유리상자"
$ws.Range("D51").Value = "에스지워너비"
$ws.Range("E51").Value = "SG워너비"
$ws.Range("B46").Value = "zJfF_41gOk8"
$ws.Range("B47").Value = "xhj-xbO6Yvs"
$ws.Range("B48").Value = "xZdTzNLmCN8"
$ws.Range("B49").Value = "LneGd_itknE"
$ws.Range("B50").Value = "URdpWdfTlao"
$ws.Range("D52").Value = "김태우"
$ws.Range("B53").Value = "pBRZzsO3L3o"
$ws.Range("D53").Value = "토이"
$ws.Range("B51").Value = "swM_GL06CxM"
$ws.Range("B52").Value = "l5BgC6iwLNA"
$ws.Range("D54").Value = "윤하"
$ws.Range("B54").Value = "FZpYfZiBEaU"

# numeric columns (ids + start times) -- order does not affect the shared
# string table, grouped here by row for clarity
$ws.Range("A39").Value = 38
$ws.Range("C39").Value = 40
$ws.Range("A40").Value = 39
$ws.Range("C40").Value = 112
$ws.Range("A41").Value = 40
$ws.Range("C41").Value = 76
$ws.Range("A42").Value = 41
$ws.Range("C42").Value = 25
$ws.Range("A43").Value = 42
$ws.Range("C43").Value = 75
$ws.Range("A44").Value = 43
$ws.Range("C44").Value = 0
$ws.Range("A45").Value = 44
$ws.Range("C45").Value = 18
$ws.Range("A46").Value = 45
$ws.Range("C46").Value = 13
$ws.Range("A47").Value = 46
$ws.Range("C47").Value = 58
$ws.Range("A48").Value = 47
$ws.Range("C48").Value = 75
$ws.Range("A49").Value = 48
$ws.Range("C49").Value = 55
$ws.Range("A50").Value = 49
$ws.Range("C50").Value = 63
$ws.Range("A51").Value = 50
$ws.Range("C51").Value = 78
$ws.Range("A52").Value = 51
$ws.Range("C52").Value = 0
$ws.Range("A53").Value = 52
$ws.Range("C53").Value = 75
$ws.Range("A54").Value = 53
$ws.Range("C54").Value = 0

# --- sheet view state: zoom, normal-view zoom, top-left cell, selection ---
$excel.ActiveWindow.Zoom = 84
$ws.Range("A26").Select()
$ws.Range("C55").Select()
